# metamiss-to-do.xlsx update
# - Row 3: mark the "no data in memory" bug as fixed, with a date in the DONE column
# - Row 4: new task about logimor(varname) behaving like logimor(varname[1]), assigned by "me",
#          noted on 2018-07-03, with an ACTION note, done on 2018-09-27
# - Row 5: new task to enable metamiss, logimor(varname) [also other options], "me", noted 2018-07-03
# - Conditional formatting (red text when E column is blank) is split: rows 3-4 keep one rule,
#   row 5 gets its own separate rule
# - Selection moves to A6 (first empty row under the table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in DONE date + note that it is fixed ---
$ws.Range("D3").Value = "fixed"
$ws.Range("E3").Value = 43370   # 2018-09-27

# --- Row 4: new task row ---
# Write D4's string first so shared-string table order matches (D4's text ends up
# before A4's text in xl/sharedStrings.xml).
$ws.Range("D4").Value = "now gives error if argument varies over observations"
$ws.Range("A4").Value = "metamiss, logimor(varname) behaves like logimor(varname[1]) "
$ws.Range("B4").Value = "me"
$ws.Range("C4").Value = 43284   # 2018-07-03
$ws.Range("E4").Value = 43370   # 2018-09-27
$ws.Rows.Item(4).RowHeight = 45

# --- Row 5: new task row (not yet actioned, D5/E5 stay blank) ---
$ws.Range("A5").Value = "enable metamiss, logimor(varname) [also other options]"
$ws.Range("B5").Value = "me"
$ws.Range("C5").Value = 43284   # 2018-07-03

# --- Split the A3:A5 conditional formatting rule into A3:A4 and A5 ---
$rngAll = $ws.Range("A3:A5")
$oldRule = $rngAll.FormatConditions.Item(1)
$oldRule.Delete()

# Create A3:A4 rule first, A5 rule second
$rngA3A4 = $ws.Range("A3:A4")
$ruleA3A4 = $rngA3A4.FormatConditions.Add(2, 0, "=ISBLANK(`$E3)")

$rngA5 = $ws.Cells.Item(5, 1)
$ruleA5 = $rngA5.FormatConditions.Add(2, 0, "=ISBLANK(`$E5)")

# Apply red font formatting - set A5's rule color first, then A3:A4's, so the
# dxf entries are created in that order (matching the target dxfId assignment)
$ruleA5.Font.Color = 255
$ruleA3A4.Font.Color = 255

# Match target priorities: A3:A4 rule keeps the (bumped) old priority, A5 rule is top priority
$ruleA3A4.Priority = 14
$ruleA5.Priority = 1

# --- Move active selection to A6 ---
$ws.Range("A6").Select()
